$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.803.07"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.778.70"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.93"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.60"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "3.209.22"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "2.792.43"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.935"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").Value = "51.727.08"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "51.65"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0462"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0845"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.27"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.98%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.03"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.89%  "
$ws.Range("D46").Value = "2.065.59"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.20"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.931"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.35%  "
